# Class02Details.docx edit
# ------------------------
# After the existing "src: having all working files" bullet, add two more
# bullets to the same list (ListParagraph style, numId=3):
#     "surge ./build"
#     "webdevbootcampclass02.surge.sh"   (this run carries a lastRenderedPageBreak)
#
# The document's "_GoBack" bookmark previously sat at the end of the
# "having all working files" paragraph (the old end of the document). Since
# that is also where editing resumed, the bookmark ends up relocated to the
# end of the new last paragraph ("webdevbootcampclass02.surge.sh") once the
# new text has been typed there.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# The paragraph we're extending is the last paragraph in the document body.
$lastPara = $d.Paragraphs.Last

# Recover that paragraph's own w:p attributes (rsid*, etc.) dynamically so
# they're preserved exactly instead of being retyped by hand.
$lastParaOpenXml = $lastPara.Range.WordOpenXML
$openTagMatch = [regex]::Match($lastParaOpenXml, '<w:p\b([^>]*)>')
$attrMatches = [regex]::Matches($openTagMatch.Groups[1].Value, '(\S+)="([^"]*)"')
$pAttrs = ""
foreach ($am in $attrMatches) {
    if ($am.Groups[1].Value.StartsWith("w:")) {
        $pAttrs += " " + $am.Groups[1].Value + "=`"" + $am.Groups[2].Value + "`""
    }
}

# 1) Re-write that paragraph's contents *without* the trailing _GoBack
#    bookmark (InsertXML replaces the full contents of the range it is
#    invoked on), keeping its run/text/pPr formatting untouched.
$srcFrag = "<w:p $wNs$pAttrs>" +
           "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
           "<w:proofErr w:type='spellStart'/><w:r><w:t>src</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
           "<w:r><w:t>: having all working files</w:t></w:r>" +
           "</w:p>"
$lastPara.Range.InsertXML($srcFrag)

# 2) Insert a new list paragraph "surge ./build" right after it.
$insPoint1 = $d.Content.End
$r1 = $d.Range($insPoint1, $insPoint1)
$surgeFrag = "<w:p $wNs>" +
             "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
             "<w:r><w:t>surge ./build</w:t></w:r>" +
             "</w:p>"
$r1.InsertXML($surgeFrag)

# 3) Insert the final list paragraph "webdevbootcampclass02.surge.sh", which
#    also picks up a lastRenderedPageBreak and the relocated _GoBack bookmark.
$insPoint2 = $d.Content.End
$r2 = $d.Range($insPoint2, $insPoint2)
$urlFrag = "<w:p $wNs>" +
           "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
           "<w:r><w:lastRenderedPageBreak/><w:t>webdevbootcampclass02.surge.sh</w:t></w:r>" +
           "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
           "</w:p>"
$r2.InsertXML($urlFrag)
